$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 now receives the data that was previously in row 3
$ws.Cells.Item(2, 2).Value = 6811909  # B2
$ws.Cells.Item(2, 3).Value = "Switzerland Challenge League"  # C2
$ws.Cells.Item(2, 4).Value = 45128.60416666666  # D2
$ws.Cells.Item(2, 5).Value = "Neuchatel Xamax"  # E2
$ws.Cells.Item(2, 6).Value = "FC Schaffhausen"  # F2
$ws.Cells.Item(2, 7).Value = 3  # G2
$ws.Cells.Item(2, 8).Value = 0  # H2
$ws.Cells.Item(2, 9).Value = 2  # I2
$ws.Cells.Item(2, 10).Value = 0  # J2
$ws.Cells.Item(2, 11).Value = "H"  # K2
$ws.Cells.Item(2, 12).Value = 2.4  # L2
$ws.Cells.Item(2, 13).Value = 3.5  # M2
$ws.Cells.Item(2, 14).Value = 2.7  # N2
$ws.Cells.Item(2, 15).Value = 2.05  # O2
$ws.Cells.Item(2, 16).Value = 3.75  # P2
$ws.Cells.Item(2, 17).Value = 3.5  # Q2
$ws.Cells.Item(2, 18).Value = -0.25  # R2
$ws.Cells.Item(2, 19).Value = 1.775  # S2
$ws.Cells.Item(2, 20).Value = 2.025  # T2
$ws.Cells.Item(2, 21).Value = 2.75  # U2
$ws.Cells.Item(2, 22).Value = 1.825  # V2
$ws.Cells.Item(2, 23).Value = 1.975  # W2
$ws.Cells.Item(2, 24).Value = 1.05  # X2
$ws.Cells.Item(2, 25).Value = -1  # Y2
$ws.Cells.Item(2, 26).Value = -1  # Z2
$ws.Cells.Item(2, 27).Value = 0.7749999999999999  # AA2
$ws.Cells.Item(2, 28).Value = -1  # AB2
$ws.Cells.Item(2, 29).Value = 0.4125  # AC2
$ws.Cells.Item(2, 30).Value = -0.5  # AD2

# Row 3 now receives the data that was previously in row 2
$ws.Cells.Item(3, 2).Value = 6811743  # B3
$ws.Cells.Item(3, 3).Value = "Switzerland Challenge League"  # C3
$ws.Cells.Item(3, 4).Value = 45128.60416666666  # D3
$ws.Cells.Item(3, 5).Value = "FC Thun"  # E3
$ws.Cells.Item(3, 6).Value = "Stade Nyonnais"  # F3
$ws.Cells.Item(3, 7).Value = 1  # G3
$ws.Cells.Item(3, 8).Value = 1  # H3
$ws.Cells.Item(3, 9).Value = 0  # I3
$ws.Cells.Item(3, 10).Value = 1  # J3
$ws.Cells.Item(3, 11).Value = "D"  # K3
$ws.Cells.Item(3, 12).Value = 1.615  # L3
$ws.Cells.Item(3, 13).Value = 4  # M3
$ws.Cells.Item(3, 14).Value = 5  # N3
$ws.Cells.Item(3, 15).Value = 1.65  # O3
$ws.Cells.Item(3, 16).Value = 4.333  # P3
$ws.Cells.Item(3, 17).Value = 4.5  # Q3
$ws.Cells.Item(3, 18).Value = -0.75  # R3
$ws.Cells.Item(3, 19).Value = 1.775  # S3
$ws.Cells.Item(3, 20).Value = 2.025  # T3
$ws.Cells.Item(3, 21).Value = 3  # U3
$ws.Cells.Item(3, 22).Value = 1.85  # V3
$ws.Cells.Item(3, 23).Value = 1.95  # W3
$ws.Cells.Item(3, 24).Value = -1  # X3
$ws.Cells.Item(3, 25).Value = 3.333  # Y3
$ws.Cells.Item(3, 26).Value = -1  # Z3
$ws.Cells.Item(3, 27).Value = -1  # AA3
$ws.Cells.Item(3, 28).Value = 1.025  # AB3
$ws.Cells.Item(3, 29).Value = -1  # AC3
$ws.Cells.Item(3, 30).Value = 0.95  # AD3

# Row 57 now receives the data that was previously in row 58
$ws.Cells.Item(57, 2).Value = 6811439  # B57
$ws.Cells.Item(57, 3).Value = "Switzerland Challenge League"  # C57
$ws.Cells.Item(57, 4).Value = 45226.60416666666  # D57
$ws.Cells.Item(57, 5).Value = "Wil 1900"  # E57
$ws.Cells.Item(57, 6).Value = "FC Baden"  # F57
$ws.Cells.Item(57, 7).Value = 0  # G57
$ws.Cells.Item(57, 8).Value = 1  # H57
$ws.Cells.Item(57, 9).Value = 0  # I57
$ws.Cells.Item(57, 10).Value = 1  # J57
$ws.Cells.Item(57, 11).Value = "A"  # K57
$ws.Cells.Item(57, 12).Value = 1.45  # L57
$ws.Cells.Item(57, 13).Value = 4.333  # M57
$ws.Cells.Item(57, 14).Value = 5.5  # N57
$ws.Cells.Item(57, 15).Value = 1.45  # O57
$ws.Cells.Item(57, 16).Value = 4.5  # P57
$ws.Cells.Item(57, 17).Value = 7  # Q57
$ws.Cells.Item(57, 18).Value = -1.25  # R57
$ws.Cells.Item(57, 19).Value = 1.9  # S57
$ws.Cells.Item(57, 20).Value = 1.9  # T57
$ws.Cells.Item(57, 21).Value = 3  # U57
$ws.Cells.Item(57, 22).Value = 1.925  # V57
$ws.Cells.Item(57, 23).Value = 1.875  # W57
$ws.Cells.Item(57, 24).Value = -1  # X57
$ws.Cells.Item(57, 25).Value = -1  # Y57
$ws.Cells.Item(57, 26).Value = 6  # Z57
$ws.Cells.Item(57, 27).Value = -1  # AA57
$ws.Cells.Item(57, 28).Value = 0.8999999999999999  # AB57
$ws.Cells.Item(57, 29).Value = -1  # AC57
$ws.Cells.Item(57, 30).Value = 0.875  # AD57

# Row 58 now receives the data that was previously in row 57
$ws.Cells.Item(58, 2).Value = 6811732  # B58
$ws.Cells.Item(58, 3).Value = "Switzerland Challenge League"  # C58
$ws.Cells.Item(58, 4).Value = 45226.60416666666  # D58
$ws.Cells.Item(58, 5).Value = "Stade Nyonnais"  # E58
$ws.Cells.Item(58, 6).Value = "AC Bellinzona"  # F58
$ws.Cells.Item(58, 7).Value = 2  # G58
$ws.Cells.Item(58, 8).Value = 3  # H58
$ws.Cells.Item(58, 9).Value = 0  # I58
$ws.Cells.Item(58, 10).Value = 0  # J58
$ws.Cells.Item(58, 11).Value = "A"  # K58
$ws.Cells.Item(58, 12).Value = 1.909  # L58
$ws.Cells.Item(58, 13).Value = 3.4  # M58
$ws.Cells.Item(58, 14).Value = 3.6  # N58
$ws.Cells.Item(58, 15).Value = 2  # O58
$ws.Cells.Item(58, 16).Value = 3.5  # P58
$ws.Cells.Item(58, 17).Value = 3.8  # Q58
$ws.Cells.Item(58, 18).Value = -0.5  # R58
$ws.Cells.Item(58, 19).Value = 2  # S58
$ws.Cells.Item(58, 20).Value = 1.8  # T58
$ws.Cells.Item(58, 21).Value = 2.5  # U58
$ws.Cells.Item(58, 22).Value = 1.8  # V58
$ws.Cells.Item(58, 23).Value = 2  # W58
$ws.Cells.Item(58, 24).Value = -1  # X58
$ws.Cells.Item(58, 25).Value = -1  # Y58
$ws.Cells.Item(58, 26).Value = 2.8  # Z58
$ws.Cells.Item(58, 27).Value = -1  # AA58
$ws.Cells.Item(58, 28).Value = 0.8  # AB58
$ws.Cells.Item(58, 29).Value = 0.8  # AC58
$ws.Cells.Item(58, 30).Value = -1  # AD58

# Row 108 now receives the data that was previously in row 109
$ws.Cells.Item(108, 2).Value = 7617773  # B108
$ws.Cells.Item(108, 3).Value = "Switzerland Challenge League"  # C108
$ws.Cells.Item(108, 4).Value = 45338.67708333334  # D108
$ws.Cells.Item(108, 5).Value = "FC Thun"  # E108
$ws.Cells.Item(108, 6).Value = "Aarau"  # F108
$ws.Cells.Item(108, 7).Value = 1  # G108
$ws.Cells.Item(108, 8).Value = 0  # H108
$ws.Cells.Item(108, 9).Value = 0  # I108
$ws.Cells.Item(108, 10).Value = 0  # J108
$ws.Cells.Item(108, 11).Value = "H"  # K108
$ws.Cells.Item(108, 12).Value = 1.727  # L108
$ws.Cells.Item(108, 13).Value = 3.8  # M108
$ws.Cells.Item(108, 14).Value = 4  # N108
$ws.Cells.Item(108, 15).Value = 1.7  # O108
$ws.Cells.Item(108, 16).Value = 4.2  # P108
$ws.Cells.Item(108, 17).Value = 4.5  # Q108
$ws.Cells.Item(108, 18).Value = -0.75  # R108
$ws.Cells.Item(108, 19).Value = 1.85  # S108
$ws.Cells.Item(108, 20).Value = 1.95  # T108
$ws.Cells.Item(108, 21).Value = 3  # U108
$ws.Cells.Item(108, 22).Value = 1.9  # V108
$ws.Cells.Item(108, 23).Value = 1.9  # W108
$ws.Cells.Item(108, 24).Value = 0.7  # X108
$ws.Cells.Item(108, 25).Value = -1  # Y108
$ws.Cells.Item(108, 26).Value = -1  # Z108
$ws.Cells.Item(108, 27).Value = 0.425  # AA108
$ws.Cells.Item(108, 28).Value = -0.5  # AB108
$ws.Cells.Item(108, 29).Value = -1  # AC108
$ws.Cells.Item(108, 30).Value = 0.8999999999999999  # AD108

# Row 109 now receives the data that was previously in row 108
$ws.Cells.Item(109, 2).Value = 7617772  # B109
$ws.Cells.Item(109, 3).Value = "Switzerland Challenge League"  # C109
$ws.Cells.Item(109, 4).Value = 45338.67708333334  # D109
$ws.Cells.Item(109, 5).Value = "FC Vaduz"  # E109
$ws.Cells.Item(109, 6).Value = "FC Sion"  # F109
$ws.Cells.Item(109, 7).Value = 1  # G109
$ws.Cells.Item(109, 8).Value = 2  # H109
$ws.Cells.Item(109, 9).Value = 1  # I109
$ws.Cells.Item(109, 10).Value = 0  # J109
$ws.Cells.Item(109, 11).Value = "A"  # K109
$ws.Cells.Item(109, 12).Value = 3.75  # L109
$ws.Cells.Item(109, 13).Value = 3.6  # M109
$ws.Cells.Item(109, 14).Value = 1.833  # N109
$ws.Cells.Item(109, 15).Value = 5.25  # O109
$ws.Cells.Item(109, 16).Value = 4  # P109
$ws.Cells.Item(109, 17).Value = 1.615  # Q109
$ws.Cells.Item(109, 18).Value = 0.75  # R109
$ws.Cells.Item(109, 19).Value = 2  # S109
$ws.Cells.Item(109, 20).Value = 1.8  # T109
$ws.Cells.Item(109, 21).Value = 2.75  # U109
$ws.Cells.Item(109, 22).Value = 1.825  # V109
$ws.Cells.Item(109, 23).Value = 1.975  # W109
$ws.Cells.Item(109, 24).Value = -1  # X109
$ws.Cells.Item(109, 25).Value = -1  # Y109
$ws.Cells.Item(109, 26).Value = 0.615  # Z109
$ws.Cells.Item(109, 27).Value = -0.5  # AA109
$ws.Cells.Item(109, 28).Value = 0.4  # AB109
$ws.Cells.Item(109, 29).Value = 0.4125  # AC109
$ws.Cells.Item(109, 30).Value = -0.5  # AD109

# Row 154 now receives the data that was previously in row 155
$ws.Cells.Item(154, 2).Value = 7617842  # B154
$ws.Cells.Item(154, 3).Value = "Switzerland Challenge League"  # C154
$ws.Cells.Item(154, 4).Value = 45405.60416666666  # D154
$ws.Cells.Item(154, 5).Value = "AC Bellinzona"  # E154
$ws.Cells.Item(154, 6).Value = "FC Schaffhausen"  # F154
$ws.Cells.Item(154, 7).Value = 1  # G154
$ws.Cells.Item(154, 8).Value = 1  # H154
$ws.Cells.Item(154, 9).Value = 0  # I154
$ws.Cells.Item(154, 10).Value = 1  # J154
$ws.Cells.Item(154, 11).Value = "D"  # K154
$ws.Cells.Item(154, 12).Value = 2.6  # L154
$ws.Cells.Item(154, 13).Value = 3.2  # M154
$ws.Cells.Item(154, 14).Value = 2.6  # N154
$ws.Cells.Item(154, 15).Value = 2.875  # O154
$ws.Cells.Item(154, 16).Value = 3  # P154
$ws.Cells.Item(154, 17).Value = 2.6  # Q154
$ws.Cells.Item(154, 18).Value = 0  # R154
$ws.Cells.Item(154, 19).Value = 2  # S154
$ws.Cells.Item(154, 20).Value = 1.8  # T154
$ws.Cells.Item(154, 21).Value = 2.25  # U154
$ws.Cells.Item(154, 22).Value = 2.025  # V154
$ws.Cells.Item(154, 23).Value = 1.775  # W154
$ws.Cells.Item(154, 24).Value = -1  # X154
$ws.Cells.Item(154, 25).Value = 2  # Y154
$ws.Cells.Item(154, 26).Value = -1  # Z154
$ws.Cells.Item(154, 27).Value = 0  # AA154
$ws.Cells.Item(154, 28).Value = 0  # AB154
$ws.Cells.Item(154, 29).Value = -0.5  # AC154
$ws.Cells.Item(154, 30).Value = 0.3875  # AD154

# Row 155 now receives the data that was previously in row 154
$ws.Cells.Item(155, 2).Value = 7617808  # B155
$ws.Cells.Item(155, 3).Value = "Switzerland Challenge League"  # C155
$ws.Cells.Item(155, 4).Value = 45405.60416666666  # D155
$ws.Cells.Item(155, 5).Value = "Neuchatel Xamax"  # E155
$ws.Cells.Item(155, 6).Value = "Wil 1900"  # F155
$ws.Cells.Item(155, 7).Value = 2  # G155
$ws.Cells.Item(155, 8).Value = 2  # H155
$ws.Cells.Item(155, 9).Value = 0  # I155
$ws.Cells.Item(155, 10).Value = 1  # J155
$ws.Cells.Item(155, 11).Value = "D"  # K155
$ws.Cells.Item(155, 12).Value = 2.15  # L155
$ws.Cells.Item(155, 13).Value = 3.3  # M155
$ws.Cells.Item(155, 14).Value = 3.1  # N155
$ws.Cells.Item(155, 15).Value = 2.7  # O155
$ws.Cells.Item(155, 16).Value = 3.3  # P155
$ws.Cells.Item(155, 17).Value = 2.625  # Q155
$ws.Cells.Item(155, 18).Value = 0  # R155
$ws.Cells.Item(155, 19).Value = 1.95  # S155
$ws.Cells.Item(155, 20).Value = 1.85  # T155
$ws.Cells.Item(155, 21).Value = 2.5  # U155
$ws.Cells.Item(155, 22).Value = 1.875  # V155
$ws.Cells.Item(155, 23).Value = 1.925  # W155
$ws.Cells.Item(155, 24).Value = -1  # X155
$ws.Cells.Item(155, 25).Value = 2.3  # Y155
$ws.Cells.Item(155, 26).Value = -1  # Z155
$ws.Cells.Item(155, 27).Value = 0  # AA155
$ws.Cells.Item(155, 28).Value = 0  # AB155
$ws.Cells.Item(155, 29).Value = 0.875  # AC155
$ws.Cells.Item(155, 30).Value = -1  # AD155

# Row 165 now receives the data that was previously in row 166
$ws.Cells.Item(165, 2).Value = 7617816  # B165
$ws.Cells.Item(165, 3).Value = "Switzerland Challenge League"  # C165
$ws.Cells.Item(165, 4).Value = 45415.63541666666  # D165
$ws.Cells.Item(165, 5).Value = "FC Baden"  # E165
$ws.Cells.Item(165, 6).Value = "FC Thun"  # F165
$ws.Cells.Item(165, 7).Value = 1  # G165
$ws.Cells.Item(165, 8).Value = 5  # H165
$ws.Cells.Item(165, 9).Value = 1  # I165
$ws.Cells.Item(165, 10).Value = 2  # J165
$ws.Cells.Item(165, 11).Value = "A"  # K165
$ws.Cells.Item(165, 12).Value = 5.5  # L165
$ws.Cells.Item(165, 13).Value = 4.5  # M165
$ws.Cells.Item(165, 14).Value = 1.5  # N165
$ws.Cells.Item(165, 15).Value = 8  # O165
$ws.Cells.Item(165, 16).Value = 5.25  # P165
$ws.Cells.Item(165, 17).Value = 1.363  # Q165
$ws.Cells.Item(165, 18).Value = 1.5  # R165
$ws.Cells.Item(165, 19).Value = 1.85  # S165
$ws.Cells.Item(165, 20).Value = 1.95  # T165
$ws.Cells.Item(165, 21).Value = 3.25  # U165
$ws.Cells.Item(165, 22).Value = 2  # V165
$ws.Cells.Item(165, 23).Value = 1.8  # W165
$ws.Cells.Item(165, 24).Value = -1  # X165
$ws.Cells.Item(165, 25).Value = -1  # Y165
$ws.Cells.Item(165, 26).Value = 0.363  # Z165
$ws.Cells.Item(165, 27).Value = -1  # AA165
$ws.Cells.Item(165, 28).Value = 0.95  # AB165
$ws.Cells.Item(165, 29).Value = 1  # AC165
$ws.Cells.Item(165, 30).Value = -1  # AD165

# Row 166 now receives the data that was previously in row 165
$ws.Cells.Item(166, 2).Value = 7617813  # B166
$ws.Cells.Item(166, 3).Value = "Switzerland Challenge League"  # C166
$ws.Cells.Item(166, 4).Value = 45415.63541666666  # D166
$ws.Cells.Item(166, 5).Value = "FC Vaduz"  # E166
$ws.Cells.Item(166, 6).Value = "Aarau"  # F166
$ws.Cells.Item(166, 7).Value = 2  # G166
$ws.Cells.Item(166, 8).Value = 1  # H166
$ws.Cells.Item(166, 9).Value = 0  # I166
$ws.Cells.Item(166, 10).Value = 0  # J166
$ws.Cells.Item(166, 11).Value = "H"  # K166
$ws.Cells.Item(166, 12).Value = 2.1  # L166
$ws.Cells.Item(166, 13).Value = 3.75  # M166
$ws.Cells.Item(166, 14).Value = 3.1  # N166
$ws.Cells.Item(166, 15).Value = 1.909  # O166
$ws.Cells.Item(166, 16).Value = 3.8  # P166
$ws.Cells.Item(166, 17).Value = 3.75  # Q166
$ws.Cells.Item(166, 18).Value = -0.5  # R166
$ws.Cells.Item(166, 19).Value = 1.9  # S166
$ws.Cells.Item(166, 20).Value = 1.9  # T166
$ws.Cells.Item(166, 21).Value = 3.25  # U166
$ws.Cells.Item(166, 22).Value = 1.925  # V166
$ws.Cells.Item(166, 23).Value = 1.875  # W166
$ws.Cells.Item(166, 24).Value = 0.909  # X166
$ws.Cells.Item(166, 25).Value = -1  # Y166
$ws.Cells.Item(166, 26).Value = -1  # Z166
$ws.Cells.Item(166, 27).Value = 0.8999999999999999  # AA166
$ws.Cells.Item(166, 28).Value = -1  # AB166
$ws.Cells.Item(166, 29).Value = -0.5  # AC166
$ws.Cells.Item(166, 30).Value = 0.4375  # AD166

# Row 167 now receives the data that was previously in row 168
$ws.Cells.Item(167, 2).Value = 7617820  # B167
$ws.Cells.Item(167, 3).Value = "Switzerland Challenge League"  # C167
$ws.Cells.Item(167, 4).Value = 45421.54166666666  # D167
$ws.Cells.Item(167, 5).Value = "Wil 1900"  # E167
$ws.Cells.Item(167, 6).Value = "FC Baden"  # F167
$ws.Cells.Item(167, 7).Value = 3  # G167
$ws.Cells.Item(167, 8).Value = 1  # H167
$ws.Cells.Item(167, 9).Value = 1  # I167
$ws.Cells.Item(167, 10).Value = 1  # J167
$ws.Cells.Item(167, 11).Value = "H"  # K167
$ws.Cells.Item(167, 12).Value = 1.475  # L167
$ws.Cells.Item(167, 13).Value = 4.1  # M167
$ws.Cells.Item(167, 14).Value = 5.2  # N167
$ws.Cells.Item(167, 15).Value = 1.4  # O167
$ws.Cells.Item(167, 16).Value = 5.25  # P167
$ws.Cells.Item(167, 17).Value = 6.5  # Q167
$ws.Cells.Item(167, 18).Value = -1.25  # R167
$ws.Cells.Item(167, 19).Value = 1.8  # S167
$ws.Cells.Item(167, 20).Value = 2  # T167
$ws.Cells.Item(167, 21).Value = 3.25  # U167
$ws.Cells.Item(167, 22).Value = 1.9  # V167
$ws.Cells.Item(167, 23).Value = 1.9  # W167
$ws.Cells.Item(167, 24).Value = 0.3999999999999999  # X167
$ws.Cells.Item(167, 25).Value = -1  # Y167
$ws.Cells.Item(167, 26).Value = -1  # Z167
$ws.Cells.Item(167, 27).Value = 0.8  # AA167
$ws.Cells.Item(167, 28).Value = -1  # AB167
$ws.Cells.Item(167, 29).Value = 0.8999999999999999  # AC167
$ws.Cells.Item(167, 30).Value = -1  # AD167

# Row 168 now receives the data that was previously in row 167
$ws.Cells.Item(168, 2).Value = 7617818  # B168
$ws.Cells.Item(168, 3).Value = "Switzerland Challenge League"  # C168
$ws.Cells.Item(168, 4).Value = 45421.54166666666  # D168
$ws.Cells.Item(168, 5).Value = "Aarau"  # E168
$ws.Cells.Item(168, 6).Value = "FC Sion"  # F168
$ws.Cells.Item(168, 7).Value = 1  # G168
$ws.Cells.Item(168, 8).Value = 2  # H168
$ws.Cells.Item(168, 9).Value = 1  # I168
$ws.Cells.Item(168, 10).Value = 1  # J168
$ws.Cells.Item(168, 11).Value = "A"  # K168
$ws.Cells.Item(168, 12).Value = 4.5  # L168
$ws.Cells.Item(168, 13).Value = 3.9  # M168
$ws.Cells.Item(168, 14).Value = 1.571  # N168
$ws.Cells.Item(168, 15).Value = 6.5  # O168
$ws.Cells.Item(168, 16).Value = 5  # P168
$ws.Cells.Item(168, 17).Value = 1.363  # Q168
$ws.Cells.Item(168, 18).Value = 1.5  # R168
$ws.Cells.Item(168, 19).Value = 1.8  # S168
$ws.Cells.Item(168, 20).Value = 2  # T168
$ws.Cells.Item(168, 21).Value = 3.25  # U168
$ws.Cells.Item(168, 22).Value = 1.975  # V168
$ws.Cells.Item(168, 23).Value = 1.825  # W168
$ws.Cells.Item(168, 24).Value = -1  # X168
$ws.Cells.Item(168, 25).Value = -1  # Y168
$ws.Cells.Item(168, 26).Value = 0.363  # Z168
$ws.Cells.Item(168, 27).Value = 0.8  # AA168
$ws.Cells.Item(168, 28).Value = -1  # AB168
$ws.Cells.Item(168, 29).Value = -0.5  # AC168
$ws.Cells.Item(168, 30).Value = 0.4125  # AD168

# Row 172 now receives the data that was previously in row 175
$ws.Cells.Item(172, 2).Value = 7617821  # B172
$ws.Cells.Item(172, 3).Value = "Switzerland Challenge League"  # C172
$ws.Cells.Item(172, 4).Value = 45429.63541666666  # D172
$ws.Cells.Item(172, 5).Value = "FC Thun"  # E172
$ws.Cells.Item(172, 6).Value = "FC Vaduz"  # F172
$ws.Cells.Item(172, 7).Value = 6  # G172
$ws.Cells.Item(172, 8).Value = 3  # H172
$ws.Cells.Item(172, 9).Value = 3  # I172
$ws.Cells.Item(172, 10).Value = 0  # J172
$ws.Cells.Item(172, 11).Value = "H"  # K172
$ws.Cells.Item(172, 12).Value = 1.5  # L172
$ws.Cells.Item(172, 13).Value = 4.5  # M172
$ws.Cells.Item(172, 14).Value = 5.25  # N172
$ws.Cells.Item(172, 15).Value = 1.45  # O172
$ws.Cells.Item(172, 16).Value = 4.75  # P172
$ws.Cells.Item(172, 17).Value = 5.5  # Q172
$ws.Cells.Item(172, 18).Value = -1.25  # R172
$ws.Cells.Item(172, 19).Value = 1.95  # S172
$ws.Cells.Item(172, 20).Value = 1.85  # T172
$ws.Cells.Item(172, 21).Value = 3.5  # U172
$ws.Cells.Item(172, 22).Value = 2  # V172
$ws.Cells.Item(172, 23).Value = 1.8  # W172
$ws.Cells.Item(172, 24).Value = 0.45  # X172
$ws.Cells.Item(172, 25).Value = -1  # Y172
$ws.Cells.Item(172, 26).Value = -1  # Z172
$ws.Cells.Item(172, 27).Value = 0.95  # AA172
$ws.Cells.Item(172, 28).Value = -1  # AB172
$ws.Cells.Item(172, 29).Value = 1  # AC172
$ws.Cells.Item(172, 30).Value = -1  # AD172

# Row 173 now receives the data that was previously in row 176
$ws.Cells.Item(173, 2).Value = 7617847  # B173
$ws.Cells.Item(173, 3).Value = "Switzerland Challenge League"  # C173
$ws.Cells.Item(173, 4).Value = 45429.63541666666  # D173
$ws.Cells.Item(173, 5).Value = "FC Schaffhausen"  # E173
$ws.Cells.Item(173, 6).Value = "FC Baden"  # F173
$ws.Cells.Item(173, 7).Value = 2  # G173
$ws.Cells.Item(173, 8).Value = 2  # H173
$ws.Cells.Item(173, 9).Value = 2  # I173
$ws.Cells.Item(173, 10).Value = 1  # J173
$ws.Cells.Item(173, 11).Value = "D"  # K173
$ws.Cells.Item(173, 12).Value = 1.571  # L173
$ws.Cells.Item(173, 13).Value = 4  # M173
$ws.Cells.Item(173, 14).Value = 4.75  # N173
$ws.Cells.Item(173, 15).Value = 1.42  # O173
$ws.Cells.Item(173, 16).Value = 4.75  # P173
$ws.Cells.Item(173, 17).Value = 6.5  # Q173
$ws.Cells.Item(173, 18).Value = -1.25  # R173
$ws.Cells.Item(173, 19).Value = 1.875  # S173
$ws.Cells.Item(173, 20).Value = 1.925  # T173
$ws.Cells.Item(173, 21).Value = 3  # U173
$ws.Cells.Item(173, 22).Value = 1.775  # V173
$ws.Cells.Item(173, 23).Value = 2.025  # W173
$ws.Cells.Item(173, 24).Value = -1  # X173
$ws.Cells.Item(173, 25).Value = 3.75  # Y173
$ws.Cells.Item(173, 26).Value = -1  # Z173
$ws.Cells.Item(173, 27).Value = -1  # AA173
$ws.Cells.Item(173, 28).Value = 0.925  # AB173
$ws.Cells.Item(173, 29).Value = 0.7749999999999999  # AC173
$ws.Cells.Item(173, 30).Value = -1  # AD173

# Row 175 now receives the data that was previously in row 173
$ws.Cells.Item(175, 2).Value = 7617823  # B175
$ws.Cells.Item(175, 3).Value = "Switzerland Challenge League"  # C175
$ws.Cells.Item(175, 4).Value = 45429.63541666666  # D175
$ws.Cells.Item(175, 5).Value = "AC Bellinzona"  # E175
$ws.Cells.Item(175, 6).Value = "FC Sion"  # F175
$ws.Cells.Item(175, 7).Value = 0  # G175
$ws.Cells.Item(175, 8).Value = 2  # H175
$ws.Cells.Item(175, 9).Value = 0  # I175
$ws.Cells.Item(175, 10).Value = 1  # J175
$ws.Cells.Item(175, 11).Value = "A"  # K175
$ws.Cells.Item(175, 12).Value = 5.75  # L175
$ws.Cells.Item(175, 13).Value = 4.333  # M175
$ws.Cells.Item(175, 14).Value = 1.444  # N175
$ws.Cells.Item(175, 15).Value = 9  # O175
$ws.Cells.Item(175, 16).Value = 4.75  # P175
$ws.Cells.Item(175, 17).Value = 1.333  # Q175
$ws.Cells.Item(175, 18).Value = 1.5  # R175
$ws.Cells.Item(175, 19).Value = 1.85  # S175
$ws.Cells.Item(175, 20).Value = 1.95  # T175
$ws.Cells.Item(175, 21).Value = 3  # U175
$ws.Cells.Item(175, 22).Value = 2  # V175
$ws.Cells.Item(175, 23).Value = 1.8  # W175
$ws.Cells.Item(175, 24).Value = -1  # X175
$ws.Cells.Item(175, 25).Value = -1  # Y175
$ws.Cells.Item(175, 26).Value = 0.333  # Z175
$ws.Cells.Item(175, 27).Value = -1  # AA175
$ws.Cells.Item(175, 28).Value = 0.95  # AB175
$ws.Cells.Item(175, 29).Value = -1  # AC175
$ws.Cells.Item(175, 30).Value = 0.8  # AD175

# Row 176 now receives the data that was previously in row 172
$ws.Cells.Item(176, 2).Value = 7617824  # B176
$ws.Cells.Item(176, 3).Value = "Switzerland Challenge League"  # C176
$ws.Cells.Item(176, 4).Value = 45429.63541666666  # D176
$ws.Cells.Item(176, 5).Value = "Aarau"  # E176
$ws.Cells.Item(176, 6).Value = "Wil 1900"  # F176
$ws.Cells.Item(176, 7).Value = 0  # G176
$ws.Cells.Item(176, 8).Value = 4  # H176
$ws.Cells.Item(176, 9).Value = 0  # I176
$ws.Cells.Item(176, 10).Value = 3  # J176
$ws.Cells.Item(176, 11).Value = "A"  # K176
$ws.Cells.Item(176, 12).Value = 2.6  # L176
$ws.Cells.Item(176, 13).Value = 3.4  # M176
$ws.Cells.Item(176, 14).Value = 2.45  # N176
$ws.Cells.Item(176, 15).Value = 2.5  # O176
$ws.Cells.Item(176, 16).Value = 3.6  # P176
$ws.Cells.Item(176, 17).Value = 2.45  # Q176
$ws.Cells.Item(176, 18).Value = 0  # R176
$ws.Cells.Item(176, 19).Value = 1.95  # S176
$ws.Cells.Item(176, 20).Value = 1.85  # T176
$ws.Cells.Item(176, 21).Value = 3.25  # U176
$ws.Cells.Item(176, 22).Value = 2  # V176
$ws.Cells.Item(176, 23).Value = 1.8  # W176
$ws.Cells.Item(176, 24).Value = -1  # X176
$ws.Cells.Item(176, 25).Value = -1  # Y176
$ws.Cells.Item(176, 26).Value = 1.45  # Z176
$ws.Cells.Item(176, 27).Value = -1  # AA176
$ws.Cells.Item(176, 28).Value = 0.8500000000000001  # AB176
$ws.Cells.Item(176, 29).Value = 1  # AC176
$ws.Cells.Item(176, 30).Value = -1  # AD176

# Row 177 now receives the data that was previously in row 181
$ws.Cells.Item(177, 2).Value = 7617848  # B177
$ws.Cells.Item(177, 3).Value = "Switzerland Challenge League"  # C177
$ws.Cells.Item(177, 4).Value = 45432.38541666666  # D177
$ws.Cells.Item(177, 5).Value = "Stade Nyonnais"  # E177
$ws.Cells.Item(177, 6).Value = "Aarau"  # F177
$ws.Cells.Item(177, 7).Value = 4  # G177
$ws.Cells.Item(177, 8).Value = 3  # H177
$ws.Cells.Item(177, 9).Value = 2  # I177
$ws.Cells.Item(177, 10).Value = 1  # J177
$ws.Cells.Item(177, 11).Value = "H"  # K177
$ws.Cells.Item(177, 12).Value = 2.25  # L177
$ws.Cells.Item(177, 13).Value = 3.6  # M177
$ws.Cells.Item(177, 14).Value = 2.7  # N177
$ws.Cells.Item(177, 15).Value = 2.05  # O177
$ws.Cells.Item(177, 16).Value = 3.8  # P177
$ws.Cells.Item(177, 17).Value = 3  # Q177
$ws.Cells.Item(177, 18).Value = -0.25  # R177
$ws.Cells.Item(177, 19).Value = 1.825  # S177
$ws.Cells.Item(177, 20).Value = 1.975  # T177
$ws.Cells.Item(177, 21).Value = 3.5  # U177
$ws.Cells.Item(177, 22).Value = 1.95  # V177
$ws.Cells.Item(177, 23).Value = 1.85  # W177
$ws.Cells.Item(177, 24).Value = 1.05  # X177
$ws.Cells.Item(177, 25).Value = -1  # Y177
$ws.Cells.Item(177, 26).Value = -1  # Z177
$ws.Cells.Item(177, 27).Value = 0.825  # AA177
$ws.Cells.Item(177, 28).Value = -1  # AB177
$ws.Cells.Item(177, 29).Value = 0.95  # AC177
$ws.Cells.Item(177, 30).Value = -1  # AD177

# Row 178 now receives the data that was previously in row 180
$ws.Cells.Item(178, 2).Value = 7617825  # B178
$ws.Cells.Item(178, 3).Value = "Switzerland Challenge League"  # C178
$ws.Cells.Item(178, 4).Value = 45432.38541666666  # D178
$ws.Cells.Item(178, 5).Value = "Wil 1900"  # E178
$ws.Cells.Item(178, 6).Value = "FC Thun"  # F178
$ws.Cells.Item(178, 7).Value = 0  # G178
$ws.Cells.Item(178, 8).Value = 3  # H178
$ws.Cells.Item(178, 9).Value = 0  # I178
$ws.Cells.Item(178, 10).Value = 0  # J178
$ws.Cells.Item(178, 11).Value = "A"  # K178
$ws.Cells.Item(178, 12).Value = 4  # L178
$ws.Cells.Item(178, 13).Value = 3.75  # M178
$ws.Cells.Item(178, 14).Value = 1.727  # N178
$ws.Cells.Item(178, 15).Value = 2.1  # O178
$ws.Cells.Item(178, 16).Value = 3.9  # P178
$ws.Cells.Item(178, 17).Value = 2.9  # Q178
$ws.Cells.Item(178, 18).Value = -0.25  # R178
$ws.Cells.Item(178, 19).Value = 1.85  # S178
$ws.Cells.Item(178, 20).Value = 1.95  # T178
$ws.Cells.Item(178, 21).Value = 3.25  # U178
$ws.Cells.Item(178, 22).Value = 1.9  # V178
$ws.Cells.Item(178, 23).Value = 1.9  # W178
$ws.Cells.Item(178, 24).Value = -1  # X178
$ws.Cells.Item(178, 25).Value = -1  # Y178
$ws.Cells.Item(178, 26).Value = 1.9  # Z178
$ws.Cells.Item(178, 27).Value = -1  # AA178
$ws.Cells.Item(178, 28).Value = 0.95  # AB178
$ws.Cells.Item(178, 29).Value = -0.5  # AC178
$ws.Cells.Item(178, 30).Value = 0.45  # AD178

# Row 179 now receives the data that was previously in row 178
$ws.Cells.Item(179, 2).Value = 7617826  # B179
$ws.Cells.Item(179, 3).Value = "Switzerland Challenge League"  # C179
$ws.Cells.Item(179, 4).Value = 45432.38541666666  # D179
$ws.Cells.Item(179, 5).Value = "FC Vaduz"  # E179
$ws.Cells.Item(179, 6).Value = "AC Bellinzona"  # F179
$ws.Cells.Item(179, 7).Value = 2  # G179
$ws.Cells.Item(179, 8).Value = 2  # H179
$ws.Cells.Item(179, 9).Value = 0  # I179
$ws.Cells.Item(179, 10).Value = 1  # J179
$ws.Cells.Item(179, 11).Value = "D"  # K179
$ws.Cells.Item(179, 12).Value = 1.85  # L179
$ws.Cells.Item(179, 13).Value = 3.6  # M179
$ws.Cells.Item(179, 14).Value = 3.6  # N179
$ws.Cells.Item(179, 15).Value = 1.727  # O179
$ws.Cells.Item(179, 16).Value = 4.333  # P179
$ws.Cells.Item(179, 17).Value = 3.8  # Q179
$ws.Cells.Item(179, 18).Value = -0.75  # R179
$ws.Cells.Item(179, 19).Value = 1.95  # S179
$ws.Cells.Item(179, 20).Value = 1.85  # T179
$ws.Cells.Item(179, 21).Value = 3.25  # U179
$ws.Cells.Item(179, 22).Value = 1.85  # V179
$ws.Cells.Item(179, 23).Value = 1.95  # W179
$ws.Cells.Item(179, 24).Value = -1  # X179
$ws.Cells.Item(179, 25).Value = 3.333  # Y179
$ws.Cells.Item(179, 26).Value = -1  # Z179
$ws.Cells.Item(179, 27).Value = -1  # AA179
$ws.Cells.Item(179, 28).Value = 0.8500000000000001  # AB179
$ws.Cells.Item(179, 29).Value = 0.8500000000000001  # AC179
$ws.Cells.Item(179, 30).Value = -1  # AD179

# Row 180 now receives the data that was previously in row 177
$ws.Cells.Item(180, 2).Value = 7617827  # B180
$ws.Cells.Item(180, 3).Value = "Switzerland Challenge League"  # C180
$ws.Cells.Item(180, 4).Value = 45432.38541666666  # D180
$ws.Cells.Item(180, 5).Value = "FC Sion"  # E180
$ws.Cells.Item(180, 6).Value = "FC Schaffhausen"  # F180
$ws.Cells.Item(180, 7).Value = 3  # G180
$ws.Cells.Item(180, 8).Value = 0  # H180
$ws.Cells.Item(180, 9).Value = 1  # I180
$ws.Cells.Item(180, 10).Value = 0  # J180
$ws.Cells.Item(180, 11).Value = "H"  # K180
$ws.Cells.Item(180, 12).Value = 1.333  # L180
$ws.Cells.Item(180, 13).Value = 4.5  # M180
$ws.Cells.Item(180, 14).Value = 8  # N180
$ws.Cells.Item(180, 15).Value = 1.333  # O180
$ws.Cells.Item(180, 16).Value = 5.5  # P180
$ws.Cells.Item(180, 17).Value = 7.5  # Q180
$ws.Cells.Item(180, 18).Value = -1.5  # R180
$ws.Cells.Item(180, 19).Value = 1.875  # S180
$ws.Cells.Item(180, 20).Value = 1.925  # T180
$ws.Cells.Item(180, 21).Value = 3.25  # U180
$ws.Cells.Item(180, 22).Value = 2  # V180
$ws.Cells.Item(180, 23).Value = 1.8  # W180
$ws.Cells.Item(180, 24).Value = 0.333  # X180
$ws.Cells.Item(180, 25).Value = -1  # Y180
$ws.Cells.Item(180, 26).Value = -1  # Z180
$ws.Cells.Item(180, 27).Value = 0.875  # AA180
$ws.Cells.Item(180, 28).Value = -1  # AB180
$ws.Cells.Item(180, 29).Value = -0.5  # AC180
$ws.Cells.Item(180, 30).Value = 0.4  # AD180

# Row 181 now receives the data that was previously in row 179
$ws.Cells.Item(181, 2).Value = 7617828  # B181
$ws.Cells.Item(181, 3).Value = "Switzerland Challenge League"  # C181
$ws.Cells.Item(181, 4).Value = 45432.38541666666  # D181
$ws.Cells.Item(181, 5).Value = "FC Baden"  # E181
$ws.Cells.Item(181, 6).Value = "Neuchatel Xamax"  # F181
$ws.Cells.Item(181, 7).Value = 2  # G181
$ws.Cells.Item(181, 8).Value = 2  # H181
$ws.Cells.Item(181, 9).Value = 1  # I181
$ws.Cells.Item(181, 10).Value = 1  # J181
$ws.Cells.Item(181, 11).Value = "D"  # K181
$ws.Cells.Item(181, 12).Value = 4.333  # L181
$ws.Cells.Item(181, 13).Value = 3.8  # M181
$ws.Cells.Item(181, 14).Value = 1.666  # N181
$ws.Cells.Item(181, 15).Value = 4.5  # O181
$ws.Cells.Item(181, 16).Value = 4.2  # P181
$ws.Cells.Item(181, 17).Value = 1.615  # Q181
$ws.Cells.Item(181, 18).Value = 1  # R181
$ws.Cells.Item(181, 19).Value = 1.8  # S181
$ws.Cells.Item(181, 20).Value = 2  # T181
$ws.Cells.Item(181, 21).Value = 3.5  # U181
$ws.Cells.Item(181, 22).Value = 1.975  # V181
$ws.Cells.Item(181, 23).Value = 1.825  # W181
$ws.Cells.Item(181, 24).Value = -1  # X181
$ws.Cells.Item(181, 25).Value = 3.2  # Y181
$ws.Cells.Item(181, 26).Value = -1  # Z181
$ws.Cells.Item(181, 27).Value = 0.8  # AA181
$ws.Cells.Item(181, 28).Value = -1  # AB181
$ws.Cells.Item(181, 29).Value = 0.9750000000000001  # AC181
$ws.Cells.Item(181, 30).Value = -1  # AD181

